# "Generate Report for Handoff"
#
# The handoff run that produced this workbook finished a little later than
# the previous one, so the three timestamp cells that record "when was the
# latest handoff xliff generated" for the last source file
# (ee252498-8cfd-416b-ad96-1bff7f2f032b.md) need to be refreshed:
#
#   - Overview!G7  (Latest HO Xliff Generate Date)      -> 2016-08-26 06:40:27
#   - zh-cn!H7     (Latest Handoff Datetime, zh-cn tab)  -> 2016-08-26 06:40:22
#   - de-de!H7     (Latest Handoff Datetime, de-de tab)  -> 2016-08-26 06:40:27
#
# These are plain text values (the cells are formatted with a date/time
# display format, but stored as shared strings), so we assign plain strings
# to keep them stored as text rather than letting Excel coerce them into
# date serial numbers.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-26 06:40:27"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-26 06:40:22"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-26 06:40:27"
